# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Updates timestamps, a few payee names, and balance figures in the
# transaction log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: timestamp corrections ---
$ws.Range("A120").Value = "2026-02-20 19:32:21"
$ws.Range("A122").Value = "2026-02-21 01:33:52"
$ws.Range("A123").Value = "2026-02-21 02:07:51"
$ws.Range("A126").Value = "2026-02-21 02:29:11"
$ws.Range("A128").Value = "2026-02-21 02:27:32"
$ws.Range("A129").Value = "2026-02-21 02:28:44"
$ws.Range("A131").Value = "2026-02-21 02:27:54"
$ws.Range("A133").Value = "2026-02-21 02:14:30"
$ws.Range("A134").Value = "2026-02-21 02:27:58"
$ws.Range("A138").Value = "2026-02-21 00:18:14"
$ws.Range("A141").Value = "2026-02-21 02:22:08"
$ws.Range("A145").Value = "2026-02-21 00:18:22"
$ws.Range("A147").Value = "2026-02-21 01:24:00"
$ws.Range("A151").Value = "2026-02-21 02:33:50"
$ws.Range("A152").Value = "2026-02-21 02:26:42"
$ws.Range("A153").Value = "2026-02-20 20:08:19"

# --- Column C: payee name corrections ---
$ws.Range("C131").Value = "MARTIAL OLIVIER TEMB ELOKAN"
$ws.Range("C145").Value = "JEAN-DANIEL NGOUFACK NGUIAZONG"
$ws.Range("C152").Value = "HERMINE NOELLE NGOMB"

# --- Column D: balance corrections ---
$ws.Range("D120").Value = 2691
$ws.Range("D122").Value = 3473
$ws.Range("D123").Value = 25233
$ws.Range("D126").Value = 87314
$ws.Range("D128").Value = 462247
$ws.Range("D129").Value = 1293
$ws.Range("D131").Value = 91412
$ws.Range("D133").Value = 115232
$ws.Range("D134").Value = 164172
$ws.Range("D138").Value = 1325
$ws.Range("D141").Value = 23092
$ws.Range("D145").Value = 172675
$ws.Range("D147").Value = 43648
$ws.Range("D151").Value = 100550
$ws.Range("D152").Value = 9887
$ws.Range("D153").Value = 51840
